$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 4 trailing rows (old rows 8-11) that no longer exist in the updated report
$ws.Range("A8:P11").EntireRow.Delete()

# Rows 3-7: refreshed monitoring data (values below mirror the latest export)

# Row 3
$ws.Range('A3').Value = 'P-22/075-S00'
$ws.Range('B3').Value = '''09-12-2022'
$ws.Range('C3').Value = '27-06-2023'
$ws.Range('D3').Value = '''104001071'
$ws.Range('E3').Value = 'TÉCNICAS REUNIDAS'
$ws.Range('F3').Value = 'Nivel'
$ws.Range('G3').Value = 'V-1040010710-0003'
$ws.Range('H3').Value = '22-075-PLN-0001'
$ws.Range('I3').Value = 'QUALITY CONTROL PLAN 22-075'
$ws.Range('J3').Value = 'PPI'
$ws.Range('K3').Value = 'Sí'
$ws.Range('L3').Value = 'Com. Menores'
$ws.Range('M3').Value = '''5'
$ws.Range('N3').Value = '24-05-2024'
$ws.Range('P3').Value = '31-10-2023 Comentado Rev. 2 // 07-11-2023 Enviado Rev. 3 // 21-12-2023 Com. Menores Rev. 3 // 05-01-2024 Enviado Rev. 4 // 07-02-2024 Com. Menores Rev. 5 // 14-02-2024 Enviado Rev. 5 // 24-05-2024 Com. Menores Rev. 5'

# Row 4
$ws.Range('A4').Value = 'P-22/075-S00'
$ws.Range('B4').Value = '''09-12-2022'
$ws.Range('C4').Value = '27-06-2023'
$ws.Range('D4').Value = '''104001071'
$ws.Range('E4').Value = 'TÉCNICAS REUNIDAS'
$ws.Range('F4').Value = 'Nivel'
$ws.Range('G4').Value = 'V-1040010710-0026'
$ws.Range('H4').Value = '22-075-PRC-0016'
$ws.Range('I4').Value = 'PAINTING PROCEDURE'
$ws.Range('J4').Value = 'Pintura'
$ws.Range('K4').Value = 'Sí'
$ws.Range('L4').Value = 'Rechazado'
$ws.Range('M4').Value = '''1'
$ws.Range('N4').Value = '''06-02-2024'
$ws.Range('O4').Value = 'Se encuentra en estado de HOLD todo el pedido'
$ws.Range('P4').Value = '26-01-2024 Enviado Rev. 0 // 26-01-2024 Rechazado Rev. 1 // 06-02-2024 Com. Mayores Rev. 1 // 06-02-2024 Rechazado Rev. 1'

# Row 5
$ws.Range('A5').Value = 'P-23/028-S00'
$ws.Range('B5').Value = '28-02-2023'
$ws.Range('C5').Value = '26-10-2023'
$ws.Range('D5').Value = '''103701061'
$ws.Range('E5').Value = 'TÉCNICAS REUNIDAS'
$ws.Range('F5').Value = 'Temperatura'
$ws.Range('G5').Value = '3998_18-1037010610-00013'
$ws.Range('H5').Value = '23-028-PRC-0009'
$ws.Range('I5').Value = 'NFXP3 - PRESERVATION AND STORAGE INSTRUCTIONS - THERMOMETERS WITH WELL'
$ws.Range('J5').Value = 'Instrucciones'
$ws.Range('K5').Value = 'No'
$ws.Range('L5').Value = 'Com. Menores'
$ws.Range('M5').Value = '''4'
$ws.Range('N5').Value = '''03-07-2024'
$ws.Range('O5').Value = 'En eGesdoc el doc. Eipsa es: 23-028-DOS-0001'
$ws.Range('P5').Value = '22-02-2024 Com. Menores Rev. 2 // 26-02-2024 Enviado Rev. 3 // 06-05-2024 Com. Menores Rev. 3 // 06-05-2024 Enviado Rev. 4 // 03-07-2024 Com. Menores Rev. 4'

# Row 6
$ws.Range('A6').Value = 'P-23/036-S00'
$ws.Range('B6').Value = '28-03-2023'
$ws.Range('C6').Value = '23-11-2023'
$ws.Range('D6').Value = 'RFQ 12-99-52-1807 _REV.A'
$ws.Range('E6').Value = 'TÉCNICAS REUNIDAS'
$ws.Range('F6').Value = 'Caudal'
$ws.Range('G6').Value = '8005710911-V-0011'
$ws.Range('H6').Value = '23-036-DOS-0002'
$ws.Range('I6').Value = 'FINAL QUALITY DOSSIER'
$ws.Range('J6').Value = 'Dossier'
$ws.Range('K6').Value = 'No'
$ws.Range('L6').Value = 'Com. Menores'
$ws.Range('M6').Value = '''0'
$ws.Range('N6').Value = '14-06-2024'
$ws.Range('O6').Value = 'Este pedido esta terminado. Mientras no reclamen no vamos ha enviar nada. Entra a fecha 14/06/2024 Aceptado con Com.Menores'
$ws.Range('P6').Value = '24-07-2023 Aprobado Rev. 0 // 14-06-2024 Com. Menores Rev. 0'

# Row 7
$ws.Range('A7').Value = 'P-23/048-S00'
$ws.Range('B7').Value = '''12-05-2023'
$ws.Range('C7').Value = '''07-01-2024'
$ws.Range('D7').Value = '''104301071'
$ws.Range('E7').Value = 'TÉCNICAS REUNIDAS'
$ws.Range('F7').Value = 'Nivel'
$ws.Range('G7').Value = '5022_20-1043010710-00004'
$ws.Range('H7').Value = '23-048-DOS-0002'
$ws.Range('I7').Value = 'NFXP4 - MANUFACTURING RECORDS BOOK FOR LEVEL GAUGES'
$ws.Range('J7').Value = 'Dossier'
$ws.Range('K7').Value = 'No'
$ws.Range('L7').Value = 'Com. Menores'
$ws.Range('M7').Value = '''0'
$ws.Range('N7').Value = '30-04-2024'
$ws.Range('O7').Value = 'Vuelve a enviar dev. 20/05/24 sin realizar ningún envío'
$ws.Range('P7').Value = '25-04-2024 Enviado Rev. 0 // 30-04-2024 Com. Menores Rev. 0'
